# "Score läggs till i Excel-filen" - add a Score column next to Name.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize the existing header casing and add the new "Score" header.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Score"

# Give Anton a score.
$ws.Range("B2").Value = 1100

# Re-fit the column widths now that there's a second column.
$ws.Columns.Item(1).ColumnWidth = 20.859375
$ws.Columns.Item(2).ColumnWidth = 9.140625
